$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "PARTNER & SENIOR SOFTWARE ENGINEER - Siege Analytics, Washington, DC | January 2014 " + [char]0x2013 + " Present",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "PARTNER - Siege Analytics, Washington, DC | January 2014 " + [char]0x2013 + " Present",
    2
)

$d.Content.Find.Execute(
    "PRINCIPAL SOFTWARE ENGINEER - Clarity and Rigour, Washington, DC | 2012 " + [char]0x2013 + " 2014",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 " + [char]0x2013 + " 2014",
    2
)

$d.Content.Find.Execute(
    "DIRECTOR OF DATA PRODUCTS - Helm, Washington, DC | 2010 " + [char]0x2013 + " 2012",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 " + [char]0x2013 + " 2012",
    2
)

$d.Content.Find.Execute(
    "SENIOR SOFTWARE ENGINEER - GSD&M, Austin, TX | 2008 " + [char]0x2013 + " 2010",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SENIOR ANALYST - Myers Research, Washington, DC | 2008 " + [char]0x2013 + " 2010",
    2
)

$d.Content.Find.Execute(
    "TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 " + [char]0x2013 + " 2004",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 " + [char]0x2013 + " 2004",
    2
)
